# Updated symbol list (GitHub Actions crypto scraper refresh).
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h) -- all stored as TEXT in the
# workbook (not numbers/percentages), so every numeric-looking value is
# written with a leading apostrophe to force Excel to keep it as a literal
# string instead of auto-coercing it to a Number.
#
# Rows 14-21 also shift down by one (a new "BitForexToken" row is inserted
# at the top of that block, pushing TigerCash/LEO/BTSEToken/... down), so
# those rows get new Coin/Link values in addition to new Price/Volume data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'307.55"
$ws.Range("E2").Value = "'-4.63%"

# Row 3 - OKB
$ws.Range("D3").Value = "'40.08"
$ws.Range("E3").Value = "'-5.85%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.013"
$ws.Range("E4").Value = "'-4.56%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.07671"
$ws.Range("E5").Value = "'-5.56%"

# Row 6 - GateToken
$ws.Range("D6").Value = "'4.225"
$ws.Range("E6").Value = "'-3.13%"

# Row 7 - FTXToken
$ws.Range("D7").Value = "'1.620"
$ws.Range("E7").Value = "'-9.80%"

# Row 8 - MXToken
$ws.Range("D8").Value = "'0.8867"
$ws.Range("E8").Value = "'-6.80%"

# Row 9 - LiechtensteinCryptoassetsExchange
$ws.Range("D9").Value = "'0.1005"
$ws.Range("E9").Value = "'-9.79%"

# Row 10 - WazirX
$ws.Range("D10").Value = "'0.1734"
$ws.Range("E10").Value = "'-6.61%"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = "'0.08947"
$ws.Range("E11").Value = "'-4.36%"

# Row 12 - BitrueCoin
$ws.Range("D12").Value = "'0.04391"
$ws.Range("E12").Value = "'-4.99%"

# Row 13 - BitMartToken
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'-0.47%"

# Row 14 - now BitForexToken (was TigerCash)
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001272"
$ws.Range("E14").Value = "'-0.99%"

# Row 15 - now TigerCash (was LEO)
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.005833"
$ws.Range("E15").Value = "'-3.49%"

# Row 16 - now LEO (was BTSEToken)
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = "'3.357"
$ws.Range("E16").Value = "'-0.71%"

# Row 17 - now BTSEToken (was BitpandaEcosystemToken)
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").Value = "'2.530"
$ws.Range("E17").Value = "'0.49%"

# Row 18 - now BitpandaEcosystemToken (was MCDex)
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").Value = "'0.3361"
$ws.Range("E18").Value = "'-0.07%"

# Row 19 - now MCDex (was ProBitToken)
$ws.Range("B19").Value = 'MCDex'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D19").Value = "'7.029"
$ws.Range("E19").Value = "'-5.61%"

# Row 20 - now ProBitToken (was ZBToken)
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = "'0.1342"
$ws.Range("E20").Value = "'-2.12%"

# Row 21 - now ZBToken (was BitForexToken)
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").Value = "'0.3030"
$ws.Range("E21").Value = "'11.72%"

# Row 22 - CoinExToken
$ws.Range("D22").Value = "'0.04240"
$ws.Range("E22").Value = "'0.88%"

# Row 23 - BitKan
$ws.Range("D23").Value = "'0.001198"
$ws.Range("E23").Value = "'-4.82%"

# Row 24 - HotbitToken
$ws.Range("D24").Value = "'0.004067"
$ws.Range("E24").Value = "'-5.40%"

# Row 25 - NitroEx (price unchanged, only volume)
$ws.Range("E25").Value = "'-6.54%"

# Row 26 - UpBots (price unchanged, only volume)
$ws.Range("E26").Value = "'-0.51%"

# Row 38 - One
$ws.Range("D38").Value = "'0.02349"
$ws.Range("E38").Value = "'-9.17%"

# Row 39 - IDEX
$ws.Range("D39").Value = "'0.05155"
$ws.Range("E39").Value = "'-6.08%"

# Row 40 - KickToken
$ws.Range("D40").Value = "'0.007962"
$ws.Range("E40").Value = "'2.14%"

# Row 41 - BKEXToken
$ws.Range("D41").Value = "'0.1324"
$ws.Range("E41").Value = "'-5.11%"

# Row 42 - Dexo
$ws.Range("D42").Value = "'0.006574"
$ws.Range("E42").Value = "'-0.14%"

# Row 43 - CEJI (price unchanged, only volume)
$ws.Range("E43").Value = "'-5.79%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.007618"
$ws.Range("E44").Value = "'-12.28%"

# Row 45 - PooCoin
$ws.Range("D45").Value = "'0.3050"
$ws.Range("E45").Value = "'-11.39%"

# Row 46 - CoinLion
$ws.Range("D46").Value = "'0.00006586"
$ws.Range("E46").Value = "'-6.24%"

# Row 47 - Kangarootoken
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.43%"

# Row 48 - BOLO
$ws.Range("D48").Value = "'0.003338"
$ws.Range("E48").Value = "'-4.30%"

# Row 49 - CoinbaseStockToken (price unchanged, only volume)
$ws.Range("E49").Value = "'40.87%"

# Row 50 - CryptobidCoin
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.43%"

# Row 51 - SpecialPowerGold
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.43%"
